# Capstone Hour Tracker - add the 2023-10-12 work-session entry (row 19)
# and the two new progress notes, matching the source XML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B19: Date (10/12/2023), C19: Start Time (11:20 AM), D19: End Time (1:10 PM)
$ws.Range("B19").Value = 45211
$ws.Range("C19").Value = 0.47222222222222227

# Column D's default style has no time number format, so copy the time
# format used by the rows above (row 18) before writing the end time.
$ws.Range("D19").NumberFormat = $ws.Range("D18").NumberFormat
$ws.Range("D19").Value = 0.54861111111111105

# G19 / H19: "What I accomplished" / "What do I think I should do next session"
$ws.Range("G19").Value = "Moved a lot of state to an editor context (I think all?), I also made a loading display for the react preview, and I moved the editor to /editor instead of /demo"
$ws.Range("H19").Value = "Next is changing the react preview to mount a basic create-react-app file template instead of npm i. Along with that I should figure out how to use pnpm instead of npm"

# Match the author's final selection in the saved workbook.
$ws.Range("H20").Select()
